$wb = $excel.ActiveWorkbook

function Edit-BomSheet {
    param($ws, [int]$jtagRow)

    $newRow = $jtagRow + 1

    # Insert a new row right below the JTAG connector row; this shifts
    # every subsequent row (footprints / unpopulated section / spacer
    # rows) down by one, matching the target layout.
    $ws.Rows($newRow).Insert()

    # --- Update the JTAG-connector row itself: switch to the SMD part ---
    $ws.Cells.Item($jtagRow, 2).WrapText = $false
    $ws.Cells.Item($jtagRow, 2).Value = "`t`n3020-10-0300-00"
    $ws.Cells.Item($jtagRow, 3).Value = 1
    $ws.Cells.Item($jtagRow, 4).Formula = '=$C$1*C' + $jtagRow
    $ws.Cells.Item($jtagRow, 5).Value = "SMD part"
    $ws.Rows($jtagRow).AutoFit()

    # --- Populate the newly inserted row with the THT alternative ---
    $ws.Cells.Item($newRow, 2).Value = "30310-6002HB"
    $ws.Cells.Item($newRow, 5).Value = "THT part alternative"
    $ws.Cells.Item($newRow, 8).Value = "JTAG Connector"
}

$ws1 = $wb.Worksheets.Item("MaxV Setup")
Edit-BomSheet $ws1 37
$ws1.Range("D18:D27").Formula = '=$C$1*C18'

$ws2 = $wb.Worksheets.Item("MaxII Setup")
Edit-BomSheet $ws2 36
$ws2.Range("D17:D26").Formula = '=$C$1*C17'
